# issue #5: stock data from json to db
#
# On the 股票 (stock) sheet, insert a new "category" column right after
# "property_category" and append two new trailing columns, "source_file"
# and "index":
#
#   before: name owner quantity face_value currency total property_category date            legislator_name legislator_id
#   after : name owner quantity face_value currency total property_category category date    legislator_name legislator_id source_file index
#
# i.e. columns I/J/K (date/legislator_name/legislator_id) shift one
# column to the right (to J/K/L), column I gets the new "category" value
# ("normal"), and two brand-new trailing columns M (source_file =
# "tmp50641") and N (index = same value as column A) are appended.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

$firstDataRow = 2
$lastDataRow = 10

# ---------------------------------------------------------------------
# Header row (row 1)
# ---------------------------------------------------------------------
# Give the three brand-new header cells (L1, M1, N1) the same look as
# the existing header cells (bold font + border), then set their text.
$ws.Range("K1").Copy($ws.Range("L1"))
$ws.Range("K1").Copy($ws.Range("M1"))
$ws.Range("K1").Copy($ws.Range("N1"))

$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# I1's text changes from "date" to "category" (it keeps its own style).
$ws.Range("I1").Value = "category"

# ---------------------------------------------------------------------
# Data rows (rows 2-10)
# ---------------------------------------------------------------------
for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    # Shift the existing date / legislator_name / legislator_id values
    # one column to the right using Copy (preserves both the original
    # value *and* its type/format, so the date text doesn't get
    # reinterpreted as a real date). Work right-to-left so a column
    # isn't overwritten before it has been copied onward.
    $ws.Cells.Item($r, 11).Copy($ws.Cells.Item($r, 12))   # K legislator_id   -> L
    $ws.Cells.Item($r, 10).Copy($ws.Cells.Item($r, 11))   # J legislator_name -> K
    $ws.Cells.Item($r, 9).Copy($ws.Cells.Item($r, 10))    # I date            -> J

    # New "category" value occupies the now-freed column I.
    $ws.Cells.Item($r, 9).Value = "normal"

    # New trailing columns: M = source_file, N = index (same as column A).
    $ws.Cells.Item($r, 1).Copy($ws.Cells.Item($r, 14))    # A index -> N
    $ws.Cells.Item($r, 13).Value = "tmp50641"
}
